$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A5 value, then copy A4 formatting (bold, centered, bordered) onto it
$ws.Range("A5").Value = "2021年"
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("B5").Value = -18.7
$ws.Range("C5").Value = 26
$ws.Range("D5").Value = -22.4
$ws.Range("E5").Value = 39.8
$ws.Range("F5").Value = -9.800000000000001
$ws.Range("G5").Value = "'"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = -23.2
$ws.Range("J5").Value = -20.2
$ws.Range("K5").Value = -29.2
$ws.Range("L5").Value = -42.9
$ws.Range("M5").Value = 22.8
$ws.Range("N5").Value = -27.9
$ws.Range("O5").Value = -10.7
$ws.Range("P5").Value = -6.7
$ws.Range("Q5").Value = 8.199999999999999
$ws.Range("R5").Value = 101.8
$ws.Range("S5").Value = 207.9
$ws.Range("T5").Value = -22.7
$ws.Range("U5").Value = -1.4
$ws.Range("V5").Value = -1
$ws.Range("W5").Value = 8.800000000000001
$ws.Range("X5").Value = 13.2
$ws.Range("Y5").Value = 8.199999999999999
$ws.Range("Z5").Value = 20.4
$ws.Range("AA5").Value = 13.9
$ws.Range("AB5").Value = 9.4
$ws.Range("AC5").Value = 10.3
$ws.Range("AD5").Value = 13.3
$ws.Range("AE5").Value = 22
$ws.Range("AF5").Value = -20.5
$ws.Range("AG5").Value = "'"
$ws.Range("AG5").Style = "Normal"
$ws.Range("AH5").Value = -56.3
$ws.Range("AI5").Value = 19.6
$ws.Range("AJ5").Value = -50
$ws.Range("AK5").Value = -39.3
$ws.Range("AL5").Value = 13.6
$ws.Range("AM5").Value = 14.2
$ws.Range("AN5").Value = -5.6
$ws.Range("AO5").Value = -5.6
$ws.Range("AP5").Value = -35.7
$ws.Range("AQ5").Value = -6.1
$ws.Range("AR5").Value = 17.1
$ws.Range("AS5").Value = -54.9
$ws.Range("AT5").Value = 160.4
$ws.Range("AU5").Value = 237.1
$ws.Range("AV5").Value = -20.5
$ws.Range("AW5").Value = -18.3
$ws.Range("AX5").Value = -20.4
$ws.Range("AY5").Value = -16
$ws.Range("AZ5").Value = 2.8
$ws.Range("BA5").Value = 1.5
$ws.Range("BB5").Value = 0.7
$ws.Range("BC5").Value = 12.9
$ws.Range("BD5").Value = 50.9
$ws.Range("BE5").Value = 8.800000000000001
$ws.Range("BF5").Value = -23.2
$ws.Range("BG5").Value = -1.4
$ws.Range("BH5").Value = -13
$ws.Range("BI5").Value = -37.5
$ws.Range("BJ5").Value = 2.1
$ws.Range("BK5").Value = -28.1
$ws.Range("BL5").Value = -13.8
$ws.Range("BM5").Value = -29.8
$ws.Range("BN5").Value = -17.9
$ws.Range("BO5").Value = -6
$ws.Range("BP5").Value = -16.5
$ws.Range("BQ5").Value = -33.3
$ws.Range("BR5").Value = 30.5
$ws.Range("BS5").Value = -7.5
$ws.Range("BT5").Value = -15.5
$ws.Range("BU5").Value = 25.1
$ws.Range("BV5").Value = -14.7
$ws.Range("BW5").Value = -13.9
$ws.Range("BX5").Value = 28.2
$ws.Range("BY5").Value = 15.7
$ws.Range("BZ5").Value = 19.5
$ws.Range("CA5").Value = -10.6
$ws.Range("CB5").Value = -85.7
$ws.Range("CC5").Value = 84.5
$ws.Range("CD5").Value = 72.5
$ws.Range("CE5").Value = 0.8
$ws.Range("CF5").Value = 9.4
$ws.Range("CG5").Value = -3
$ws.Range("CH5").Value = -4.3
$ws.Range("CI5").Value = 20.4
$ws.Range("CJ5").Value = -49.6
$ws.Range("CK5").Value = 1.3
$ws.Range("CL5").Value = 7.7
$ws.Range("CM5").Value = -26.5
$ws.Range("CN5").Value = 47
$ws.Range("CO5").Value = 11.5
$ws.Range("CP5").Value = 37.6
$ws.Range("CQ5").Value = -14.9
$ws.Range("CR5").Value = -23.9
$ws.Range("CS5").Value = 8.699999999999999
$ws.Range("CT5").Value = 9.6
$ws.Range("CU5").Value = 23.2
$ws.Range("CV5").Value = -13.8
$ws.Range("CW5").Value = -1.6
$ws.Range("CX5").Value = 22.7
$ws.Range("CY5").Value = -39.8
$ws.Range("CZ5").Value = 30.3
$ws.Range("DA5").Value = 5.3
$ws.Range("DB5").Value = 4.8
$ws.Range("DC5").Value = -7.4
$ws.Range("DD5").Value = -15.6
$ws.Range("DE5").Value = -11.1
$ws.Range("DF5").Value = 10.6
$ws.Range("DG5").Value = 25.7
$ws.Range("DH5").Value = 7.3
$ws.Range("DI5").Value = -6.7
$ws.Range("DJ5").Value = 20.6
$ws.Range("DK5").Value = 46.4
